# The table originally occupied columns B:F (with a redundant helper
# value duplicated in column A). The edit removes column A entirely,
# shifting every remaining column (B->A, C->B, D->C, E->D, F->E) one
# position to the left, so the real header ("QS_Astral_exact50") lands
# in A1 and the data lines up in A:E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()
